# update W6D1 Pandas notebook
# Sheet1 holds a small DataFrame export (columns B:E = Id, Python,
# Machine Learning, Deep Learning; rows 2-7 are the data rows).
# The edit fills in the previously-missing Deep Learning score for
# row 4 (E4 = 5) and clears the stray value that had been left in
# row 6 (E6), leaving that Deep Learning cell blank again.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = 5
$ws.Range("E6").ClearContents()
